$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns at the front (A, B) - shifts nome/modelo/preco/... right by 2
$ws.Columns("A:B").Insert()

# 2. New headers for the inserted columns
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"

# Copy header style (bold + border) from the neighbouring header cell onto the new headers
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# 3. Fill the new "data"/"loja" columns for every existing data row (rows 2-6 at this point)
$ws.Range("A2:A6").Value = "30/07/2024"
$ws.Range("B2:B6").Value = "subgrave autoparts"

# 4. Remove the rows that are no longer present in the sheet
#    (the "K1200" row and the "Voltimetro Sequenciador Digital" row)
$ws.Rows("5").Delete()
$ws.Rows("4").Delete()

# F4 ("politica") stays blank for this row, same as it was before the edit
$ws.Range("F4").ClearContents()

# 5. Update the remaining rows' link column (now column I) with the refreshed URLs
$ws.Range("I2").Value = "https://produto.mercadolivre.com.br/MLB-2873301438-controle-longa-distancia-jfa-redline-wr-key1-ad1-multimidia-_JM#position%3D46%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D2ea3a4f5-9fab-41e6-9455-ea1bea1f7731"
$ws.Range("I3").Value = "https://produto.mercadolivre.com.br/MLB-2872406163-controle-longa-distancia-jfa-redline-wr-p-aparelho-original-_JM#position%3D20%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Dfac91a18-86de-4e6f-b4f0-68dd328fa5ee"
$ws.Range("I4").Value = "https://produto.mercadolivre.com.br/MLB-3185900332-voltimetro-jfa-vs5hi-3-em-1-sequenciador-high-voltagem-12v-_JM#position%3D31%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D0d5ca327-689d-4181-ac63-44ca774b4075"
